$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the existing Key Message text for row 2 (B2)
$ws.Range("B2").Value = "Kids are  celebrating Holi ,also there were poor children looking at them,  Be kind to poor people"

# Add a new row of data (row 3): S.No. = 2, Theme = new shared string text
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Children enjoying on village side , they were plaing games on river bed"

# Adjust row height for row 2 from 45 to 30
$ws.Rows.Item(2).RowHeight = 30

# Update the selected cell/range
$ws.Range("B3").Select()
